# Adding Group 25 untestables
#
# For each requirement row listed below, mark column F ("Unable to Test")
# the same way it is already marked for the existing I/K "Unable To Test"
# columns elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5,6,7,11,12,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,31,32,33,35,36,37,38,39,41,42,43,46,47,48,49,50,51,52,58,60,62,63,64,65,66,69,70,71,72,73,74,75,77,78,79,80,81,82,83,86,87,88,90,93,94,97,101,102,106,109,110,113,116,120,121,122,123,124,125,126,127,128,129,131,132,133,134,137,138,140,141,143,148,149,152,153,156,157,158,159,160,161,163,164,165,166,167,168,169,170,171,172,173,174,179,182,183,184,185,186,187,188,189,190,191,192,193,194,196,197,198,199,200,201,202,203,205,206,207,210,216,217,218)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "Unable to Test"
}

# Restore the view state captured in the saved workbook: zoomed in on the
# sheet, scrolled so column B is left-most, and selected near the bottom of
# the requirements list.
$excel.ActiveWindow.Zoom = 180
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F225").Select()
